# Applies the "extract MINICOG/SLUMS/BNT/BVRT" edit:
#  - Renames header A1 from "NewOverallID" to "NoteID"
#  - Replaces the separate PHYSICIAN PROBLEM LIST / Mental status rows
#    (rows 17-27) with a single combined multi-line "note" cell in B17
#    that wraps text, and clears the now-unused rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("A1").Value = "NoteID"

# Build the combined note text (uses vbLf between the original separate lines)
$noteLines = @(
    "PHYSICIAN PROBLEM LIST:",
    "MMSE 28/30 (October 2006 -- 1/3 spointaneous recall); clock drawing intact, 25/30 (June 2008)",
    "B12 540",
    "TSH 1.45",
    "MRI: small vessel ischemic changes",
    "REPEAT evaluation 11/2007: MMSE 24/30...",
    [char]0x2026,
    "Mental Status Exam:",
    "...",
    "MMSE 25/30",
    [char]0x2026
)
$noteText = [string]::Join([char]10, $noteLines)

# Row 17 becomes the combined note, with wrap text + taller row
$ws.Cells.Item(17, 2).Value = $noteText
$ws.Rows.Item(17).RowHeight = 191.25
$ws.Cells.Item(17, 2).WrapText = $true
$ws.Cells.Item(17, 2).VerticalAlignment = -4108  # xlCenter

# Row 18 now holds the next record id (17) in column A, and an empty B cell
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).ClearContents()

# Rows 19-27: clear column B (previously held separate note fragments)
for ($r = 19; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).ClearContents()
}

# Update selection / view to match the saved workbook state
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 15
